# Apply updated crypto market data values (prices & 1h volume change)
# to match the refreshed GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.967.19"
$ws.Range("E2").Value = "  +3.12%  "
$ws.Range("D3").Value = "2.454.56"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'578.07"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "2.453.51"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D12").Value = "'5.29"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D14").Value = "'28.47"
$ws.Range("E14").Value = "  +7.64%  "
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("D17").Value = "62.905.55"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").Value = "2.453.26"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "'7.97"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "'11.11"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("D21").Value = "'332.00"
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").Value = "'2.11"
$ws.Range("E23").Value = "  +10.67%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'66.35"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.25"
$ws.Range("E26").Value = "  +25.03%  "
$ws.Range("D27").Value = "'648.61"
$ws.Range("E27").Value = "  +10.16%  "
$ws.Range("E28").Value = "  +4.39%  "
$ws.Range("E29").Value = "  +5.99%  "
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").Value = "'1.44"
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("D35").Value = "0.0₆0414"
$ws.Range("E35").Value = "  +46.56%  "
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D38").Value = "'4.78"
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("D41").Value = "'152.47"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'18.85"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = "  +10.23%  "
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("D45").Value = "'42.70"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'14.99"
$ws.Range("E47").Value = "  +27.34%  "
$ws.Range("D48").Value = "'146.33"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").Value = "'20.74"
$ws.Range("E50").Value = "  +5.37%  "
